# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" right before the existing "2022-Q1"
#   sheet (so tab order becomes 总计, 2022-Q3, 2022-Q1, 2021-Q4, 2021-Q3)
#   and populate it with the fund-holdings table for that quarter.
# - Insert a new row into the "总计" (totals) sheet for 2022-Q3 and
#   renumber the existing rows' index column accordingly.

function Set-TextValue {
    # Force a value to be stored as text (keeps leading zeros / trailing
    # zeros like "011531" or "59.37" intact instead of Excel coercing it
    # to a number), then strip the "Text" number-format Excel auto-applies
    # so the cell keeps the workbook's default (unstyled) look.
    param($range, $value)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert the 2022-Q3 row at the top of the data and
#    renumber the following rows.
# ---------------------------------------------------------------------
$total.Rows(2).Insert()

# New row 2 picks up column-A's index style (bold/centered/bordered,
# same style used by the other id cells) by copying format from the
# (now shifted-down) old row-2 id cell.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)
$total.Cells.Item(2, 1).Value = 0

# The freshly inserted row also inherited a blended style on B:D from
# the row above (the header row) - reset that back to the plain/default
# look used by every other data row.
$total.Range("B2:D2").Style = "Normal"
Set-TextValue $total.Range("B2") "2022-Q3"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 3.44

# Renumber the id column for the rows that shifted down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3

# ---------------------------------------------------------------------
# 2. Add the new "2022-Q3" sheet (placed right before "2022-Q1") and
#    fill in the fund table.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q3 = $wb.Worksheets.Add($q1)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    Set-TextValue $q3.Cells.Item(1, $col) $headers[$col - 2]
}
# Copy the bold/centered/bordered header style from the 总计 sheet onto
# the whole header row (styles are shared workbook-wide).
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$rows = @(
    @("011531", "朱雀恒心一年持有期混合",                         "59.37", "93.18", "2.38", "1.4130", 10),
    @("007493", "朱雀产业臻选混合A",                               "34.47", "92.95", "2.57", "0.8859", 10),
    @("010922", "朱雀匠心一年持有期混合",                         "15.23", "91.81", "3.00", "0.4569", 8),
    @("007494", "朱雀产业臻选混合C",                               "11.69", "92.95", "2.57", "0.3004", 10),
    @("008294", "朱雀企业优胜股票A",                               "11.46", "93.76", "2.60", "0.2980", 10),
    @("008295", "朱雀企业优胜股票C",                               "2.16",  "93.76", "2.60", "0.0562", 10),
    @("004266", "招商沪港深科技创新主题精选灵活配置混合A", "0.92",  "90.52", "2.26", "0.0208", 10),
    @("010754", "招商沪港深科技创新主题精选灵活配置混合C", "0.25",  "90.52", "2.26", "0.0056", 10)
)

$r = 2
foreach ($row in $rows) {
    Set-TextValue $q3.Cells.Item($r, 2) $row[0]
    Set-TextValue $q3.Cells.Item($r, 3) $row[1]
    Set-TextValue $q3.Cells.Item($r, 4) $row[2]
    Set-TextValue $q3.Cells.Item($r, 5) $row[3]
    Set-TextValue $q3.Cells.Item($r, 6) $row[4]
    Set-TextValue $q3.Cells.Item($r, 7) $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Id column (A2:A9): 0-based index, same style as the 总计 sheet's id
# column.
for ($i = 2; $i -le 9; $i++) {
    $q3.Cells.Item($i, 1).Value = $i - 2
}
$q3.Range("A2:A9").Style = "Normal"
$total.Range("A2").Copy()
$q3.Range("A2:A9").PasteSpecial(-4122)

# Restore the originally-active tab (2021-Q3) - adding/renaming sheets
# shifts Excel's notion of the "active" sheet to whatever was touched last.
$wb.Worksheets.Item("2021-Q3").Activate()
